$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Moisture Units" row (row 24: ~moisturetunits / From Roast>Properties>Moisture / pct)
$ws.Rows.Item(24).Delete()

# Scroll the view down so row 2 is at the top, with A2 selected (not the A2:C2 range)
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
